$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 342 (shifts existing rows 342-391 down to 343-392,
# and keeps formatting inherited from the surrounding rows).
$ws.Rows.Item(342).EntireRow.Insert()

# Populate the new row 342 with the new daily price record.
$ws.Range("A342").Value2 = 5
$ws.Range("B342").Value = "Macroferia Regional de Talca"
$ws.Range("C342").Value = "Maule"
$ws.Range("D342").Value2 = 44491
$ws.Range("E342").Value2 = 7
$ws.Range("F342").Value = "Fruta"
$ws.Range("G342").Value2 = 100108
$ws.Range("H342").Value = "Tropicales y subtropicales"
$ws.Range("I342").Value2 = 100108006
$ws.Range("J342").Value = "Plátano"
$ws.Range("K342").Value = "Sin especificar"
$ws.Range("L342").Value = "Pintón"
$ws.Range("M342").Value2 = 1140
$ws.Range("N342").Value2 = 22000
$ws.Range("O342").Value2 = 23000
$ws.Range("P342").Value2 = 22526
$ws.Range("Q342").Value = "$/caja 20 kilos"
$ws.Range("R342").Value = "Ecuador"
$ws.Range("S342").Value2 = 1126
$ws.Range("T342").Value2 = 20
